$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "NSE:LICMFGOLD"
$ws.Range("C2").Value = "NSE:AUROPHARMA"
$ws.Range("E2").Value = "NSE:ABCAPITAL"
$ws.Range("F2").Value = ""

# Row 3
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = "NSE:BAJAJHCARE"
$ws.Range("E3").Value = "NSE:AMBUJACEM"

# Row 4
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = "NSE:CHALET"
$ws.Range("E4").Value = "NSE:ANGELONE"

# Row 5
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = "NSE:COASTCORP"
$ws.Range("E5").Value = "NSE:APOLLOHOSP"

# Row 6
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = "NSE:ESTER"
$ws.Range("E6").Value = "NSE:ASTRAL"

# Row 7
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = "NSE:JAICORPLTD"
$ws.Range("E7").Value = "NSE:BOSCHLTD"

# Row 8
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = "NSE:KANSAINER"
$ws.Range("E8").Value = "NSE:BPCL"

# Row 9
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = "NSE:KEYFINSERV"
$ws.Range("E9").Value = "NSE:CONCOR"

# Row 10
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = "NSE:KPRMILL"
$ws.Range("E10").Value = "NSE:CYIENT"

# Row 11
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = "NSE:LINCOLN"
$ws.Range("E11").Value = "NSE:DEEPAKNTR"

# Row 12
$ws.Range("C12").Value = "NSE:MONARCH"
$ws.Range("E12").Value = "NSE:DRREDDY"

# Row 13
$ws.Range("C13").Value = "NSE:MVGJL"
$ws.Range("E13").Value = "NSE:EICHERMOT"

# Row 14
$ws.Range("C14").Value = "NSE:NDL"
$ws.Range("E14").Value = "NSE:HAL"

# Row 15
$ws.Range("C15").Value = "NSE:OILCOUNTUB"
$ws.Range("E15").Value = "NSE:IGL"

# Row 16
$ws.Range("C16").Value = "NSE:ONGC"
$ws.Range("E16").Value = "NSE:PETRONET"

# Row 17
$ws.Range("C17").Value = "NSE:ONWARDTEC"
$ws.Range("E17").Value = "NSE:POLYCAB"

# Row 18 and 19 (new rows) - copy the style of A17 (bordered/bold/centered)
# down into A18:A19 so the new index cells match the existing column-A style
$ws.Range("A17").Copy()
$ws.Range("A18:A19").PasteSpecial(-4122)

# Row 18 (new)
$ws.Range("A18").Value = 16
$ws.Range("C18").Value = "NSE:PVRINOX"
$ws.Range("E18").Value = "NSE:POONAWALLA"

# Row 19 (new)
$ws.Range("A19").Value = 17
$ws.Range("C19").Value = "NSE:RACE"
$ws.Range("E19").Value = "NSE:POWERGRID"
